# Finish descriptive feature name column in docs
#
# Fills in the previously-empty "Type" / "Descriptive Name" / "Nulls"
# cells of the single data-dictionary table, per the commit diff.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $text) {
    $cell = $table.Cell($row, $col)
    $cell.Range.Text = $text
}

# --- Row F0 : Datapoint Number / Row Number ---------------------------
# Type: "Int [0,]" -> "Int [0,+]"
Set-CellText $t 2 4 "Int [0,+]"

# --- Row F3 : PM25 ------------------------------------------------------
# Type: "[" -> "Float [0,+]"
Set-CellText $t 5 4 "Float [0,+]"

# --- Row F4 : Date -------------------------------------------------------
# Type: "" -> "Datetime"
Set-CellText $t 6 4 "Datetime"

# --- Row F5 : temp ---------------------------------------------------------
# Type: "" -> "Float [-,100]"
Set-CellText $t 7 4 "Float [-,100]"

# --- Row F6 : dewpoint ------------------------------------------------------
# Type: "" -> "Float [-,+]"
Set-CellText $t 8 4 "Float [-,+]"

# --- Row F7 : RH -------------------------------------------------------------
# Type: "" -> "[0,100]"
Set-CellText $t 9 4 "[0,100]"

# --- Row F32 : Clouds / Cloud Coverage ---------------------------------------
# Type: "" -> "[0,8]"; Nulls: "" -> "Yes"
Set-CellText $t 34 4 "[0,8]"
Set-CellText $t 34 5 "Yes"

# --- Rows F33..F41 : Clds1000 .. Clds10000 -----------------------------------
# Descriptive Name: "" -> "Cloud Coverage at N000 ft"; Type -> "[0,8]"; Nulls -> "Yes"
$cloudRows = @(
    @{ Row = 35; Alt = "1000" },
    @{ Row = 36; Alt = "2000" },
    @{ Row = 37; Alt = "3000" },
    @{ Row = 38; Alt = "4000" },
    @{ Row = 39; Alt = "5000" },
    @{ Row = 40; Alt = "6000" },
    @{ Row = 41; Alt = "7000" },
    @{ Row = 42; Alt = "8000" },
    @{ Row = 43; Alt = "9000" },
    @{ Row = 44; Alt = "10000" }
)

foreach ($entry in $cloudRows) {
    $rowIdx = $entry.Row
    $alt = $entry.Alt
    Set-CellText $t $rowIdx 3 "Cloud Coverage at $alt ft"
    Set-CellText $t $rowIdx 4 "[0,8]"
    Set-CellText $t $rowIdx 5 "Yes"
}

Write-Host "Data dictionary table filled in."
